$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are numeric-looking strings that must stay as text.
# Force text format, assign, then restore default style so no stray formatting remains.
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "29.156.95"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.824.50"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9990"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "234.79"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6007"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.2790"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "23.49"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.07603"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.825.35"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.789"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6293"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.000009902"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "78.81"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.850"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "29.162.93"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "226.10"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.71"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.991"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "154.90"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.015"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.1299"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "16.55"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.488"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.06254"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.449"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.829"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.795"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.121"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.737"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.6389"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.531"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.215.96"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.723"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.01731"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.495"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9033"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.9992"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.984.90"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "100.29"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "62.72"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.00000000115"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.596"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.492"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.4549"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.05499"
$cell.Style = "Normal"

# Columns B, C, E are plain text already and do not need coercion guards.
$ws.Cells.Item(2, 5).Value = "  -0.58%  "
$ws.Cells.Item(3, 5).Value = "  -0.78%  "
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(6, 5).Value = "  -3.95%  "
$ws.Cells.Item(7, 5).Value = "  -0.02%  "
$ws.Cells.Item(8, 5).Value = "  -4.75%  "
$ws.Cells.Item(9, 5).Value = "  -3.46%  "
$ws.Cells.Item(10, 5).Value = "  -5.06%  "
$ws.Cells.Item(12, 5).Value = "  -0.75%  "
$ws.Cells.Item(13, 5).Value = "  -3.19%  "
$ws.Cells.Item(14, 5).Value = "  -6.48%  "
$ws.Cells.Item(15, 5).Value = "  -2.64%  "
$ws.Cells.Item(17, 5).Value = "  -3.42%  "
$ws.Cells.Item(18, 5).Value = "  -5.73%  "
$ws.Cells.Item(19, 5).Value = "  -0.64%  "
$ws.Cells.Item(20, 5).Value = "  -2.39%  "
$ws.Cells.Item(21, 5).Value = "  -0.02%  "
$ws.Cells.Item(22, 5).Value = "  -4.57%  "
$ws.Cells.Item(23, 5).Value = "  -4.60%  "
$ws.Cells.Item(24, 5).Value = "  -0.07%  "
$ws.Cells.Item(25, 5).Value = "  -2.00%  "
$ws.Cells.Item(26, 5).Value = "  -5.22%  "
$ws.Cells.Item(27, 5).Value = "  -3.32%  "
$ws.Cells.Item(28, 5).Value = "  -4.43%  "
$ws.Cells.Item(29, 5).Value = "  +2.16%  "
$ws.Cells.Item(30, 5).Value = "  -13.76%  "
$ws.Cells.Item(31, 5).Value = "  -1.94%  "
$ws.Cells.Item(32, 5).Value = "  -5.09%  "
$ws.Cells.Item(33, 5).Value = "  -5.94%  "
$ws.Cells.Item(34, 5).Value = "  -1.50%  "
$ws.Cells.Item(35, 5).Value = "  -4.26%  "
$ws.Cells.Item(36, 5).Value = "  -7.95%  "
$ws.Cells.Item(37, 5).Value = "  -1.48%  "
$ws.Cells.Item(38, 5).Value = "  -1.00%  "
$ws.Cells.Item(39, 5).Value = "  -3.11%  "
$ws.Cells.Item(40, 5).Value = "  -5.39%  "
$ws.Cells.Item(41, 5).Value = "  -6.21%  "
$ws.Cells.Item(42, 5).Value = "  -4.18%  "
$ws.Cells.Item(43, 5).Value = "  -0.07%  "
$ws.Cells.Item(44, 5).Value = "  -0.28%  "
$ws.Cells.Item(45, 5).Value = "  -0.06%  "
$ws.Cells.Item(46, 5).Value = "  -3.96%  "
$ws.Cells.Item(47, 5).Value = "  -3.46%  "
$ws.Cells.Item(48, 2).Value = "RenderToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(48, 5).Value = "  -6.07%  "
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 5).Value = "  -4.33%  "
$ws.Cells.Item(50, 5).Value = "  -0.76%  "
$ws.Cells.Item(51, 5).Value = "  -2.64%  "
